# The paragraph "<id>p026r_1</id>" was previously split across five runs
# (<id> / p026 / r / _1 / </id>), each with its own formatting. Collapse
# them into a single run holding the full text, using the same formatting
# as the original "<id>" / "</id>" runs (Courier New, color 7f6000, 9pt).
# A Find/Replace whose search text spans the run boundaries makes Word
# merge the matched runs into one, adopting the formatting of the first
# run in the match - exactly the effect we want here.
$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p026r_1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p026r_1</id>", 2)
